# Sprint 2, day 3 daily standup update.
# Update the "Sprint" worksheet with the day-3 status / effort entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# Row 5: "Create app logo" (Đuro Nenadović) moved to "In progress";
# 1 day of effort logged against Day 3 (column I).
$ws.Cells.Item(5, 6).Value = "In progress"
$ws.Cells.Item(5, 9).Value = 1

# Row 6: "Design app gui" (Vanja Cvetković) completed;
# remaining 5 days of effort logged against Day 3 (column I).
$ws.Cells.Item(6, 6).Value = "Done"
$ws.Cells.Item(6, 9).Value = 5

# Row 11: "Expand timeline for tracks" assigned to Vanja Cvetković,
# and moved to "In progress".
$ws.Cells.Item(11, 4).Value = "Vanja Cvetković"
$ws.Cells.Item(11, 6).Value = "In progress"

# Move the active selection to reflect where the user was working.
$ws.Range("I6").Select()

$excel.ActiveWorkbook.Save()
